$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G gets a "retour sur l'existant" status for every requirement row:
# functional requirements (rows 2-17) are marked "A faire" ...
$ws.Range("G2:G17").Value = "A faire"
# ...while the non-functional requirements (rows 18-26) are marked "?"
$ws.Range("G18:G26").Value = "?"

# Widen column G so the new values fit (target stored width ~27.78 characters)
$ws.Range("G1:G26").ColumnWidth = 27.06305803571428

# Move the active selection to U3, matching the author's final cursor position
[void]$ws.Range("U3").Select()
